$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Jaren Jackson Jr."
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "Memphis Grizzlies"

$ws.Range("A10").Value = "Jayson Tatum"
$ws.Range("B10").Value = "SF,PF"
$ws.Range("C10").Value = "Boston Celtics"

$ws.Range("A11").Value = "Keyonte George"
$ws.Range("B11").Value = "PG,SG"
$ws.Range("C11").Value = "Utah Jazz"

$ws.Range("A12").Value = "Ivica Zubac"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "LA Clippers"

$ws.Range("A13").Value = "Zion Williamson"
$ws.Range("B13").Value = "PF,C"
$ws.Range("C13").Value = "New Orleans Pelicans"
